# chore: update Sheets via scheduled runner
# Refresh cached market-board price / profit figures on each class's Leve
# profit sheet (currentAveragePrice*, LevePrice*, LeveProfit* columns).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H57").Value = 61000
$ws.Range("J57").Value = 61000
$ws.Range("L57").Value = 183000
$ws.Range("N57").Value = -183998
$ws.Range("H74").Value = 9059.5625
$ws.Range("I74").Value = 8180.6924
$ws.Range("J74").Value = 12868
$ws.Range("K74").Value = 8180.6924
$ws.Range("L74").Value = 12868
$ws.Range("M74").Value = -7244.6924
$ws.Range("N74").Value = -14740
$ws.Range("H77").Value = 9059.5625
$ws.Range("I77").Value = 8180.6924
$ws.Range("J77").Value = 12868
$ws.Range("K77").Value = 40903.462
$ws.Range("L77").Value = 64340
$ws.Range("M77").Value = -36223.462
$ws.Range("N77").Value = -73700
$ws.Range("H92").Value = 2145.8333
$ws.Range("I92").Value = 65.25
$ws.Range("K92").Value = 65.25
$ws.Range("M92").Value = 1182.75
$ws.Range("H94").Value = 4262.5
$ws.Range("I94").Value = 5250
$ws.Range("J94").Value = 1300
$ws.Range("K94").Value = 5250
$ws.Range("L94").Value = 1300
$ws.Range("M94").Value = -4799
$ws.Range("N94").Value = -2202
$ws.Range("H100").Value = 6239.6
$ws.Range("I100").Value = 3001.3333
$ws.Range("K100").Value = 3001.3333
$ws.Range("M100").Value = -2460.3333
$ws.Range("H112").Value = 1610.1
$ws.Range("J112").Value = 1654.0625
$ws.Range("L112").Value = 4962.1875
$ws.Range("N112").Value = -7178.1875
$ws.Range("H132").Value = 1912.8148
$ws.Range("I132").Value = 1686.6842
$ws.Range("J132").Value = 2449.875
$ws.Range("K132").Value = 5060.0526
$ws.Range("L132").Value = 7349.625
$ws.Range("M132").Value = -2530.0526
$ws.Range("N132").Value = -12409.625
$ws.Range("H138").Value = 3525.2126
$ws.Range("I138").Value = 3106.8
$ws.Range("J138").Value = 3721.3438
$ws.Range("K138").Value = 9320.400000000001
$ws.Range("L138").Value = 11164.0314
$ws.Range("M138").Value = -4180.400000000001
$ws.Range("N138").Value = -21444.0314
$ws.Range("H139").Value = 69993.625
$ws.Range("J139").Value = 69993.625
$ws.Range("L139").Value = 69993.625
$ws.Range("N139").Value = -80273.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 8133557.5
$ws.Range("I74").Value = 9012283
$ws.Range("J74").Value = 5350
$ws.Range("K74").Value = 9012283
$ws.Range("L74").Value = 5350
$ws.Range("M74").Value = -9011409
$ws.Range("N74").Value = -7098
$ws.Range("H77").Value = 8133557.5
$ws.Range("I77").Value = 9012283
$ws.Range("J77").Value = 5350
$ws.Range("K77").Value = 45061415
$ws.Range("L77").Value = 26750
$ws.Range("M77").Value = -45057047
$ws.Range("N77").Value = -35486
$ws.Range("H97").Value = 833.43335
$ws.Range("I97").Value = 868.2273
$ws.Range("J97").Value = 737.75
$ws.Range("K97").Value = 868.2273
$ws.Range("L97").Value = 737.75
$ws.Range("M97").Value = -372.2273
$ws.Range("N97").Value = -1729.75
$ws.Range("H132").Value = 3339.6
$ws.Range("I132").Value = 1422.75
$ws.Range("K132").Value = 4268.25
$ws.Range("M132").Value = -1738.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 4984
$ws.Range("I26").Value = 4984
$ws.Range("K26").Value = 4984
$ws.Range("M26").Value = -4692
$ws.Range("H81").Value = 44657.832
$ws.Range("J81").Value = 44657.832
$ws.Range("L81").Value = 44657.832
$ws.Range("N81").Value = -46779.832
$ws.Range("H84").Value = 44657.832
$ws.Range("J84").Value = 44657.832
$ws.Range("L84").Value = 133973.496
$ws.Range("N84").Value = -144581.496
$ws.Range("H94").Value = 1116.6666
$ws.Range("J94").Value = 1500
$ws.Range("L94").Value = 1500
$ws.Range("N94").Value = -2402
$ws.Range("H99").Value = 1541.5
$ws.Range("I99").Value = 1506.75
$ws.Range("K99").Value = 1506.75
$ws.Range("M99").Value = -8.75
$ws.Range("H107").Value = 4805.75
$ws.Range("I107").Value = 4924.4546
$ws.Range("K107").Value = 4924.4546
$ws.Range("M107").Value = -3004.4546
$ws.Range("H134").Value = 4101.8335
$ws.Range("I134").Value = 1871.2858
$ws.Range("K134").Value = 5613.857400000001
$ws.Range("M134").Value = -3078.857400000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 814
$ws.Range("I22").Value = 638.8148
$ws.Range("J22").Value = 1208.1666
$ws.Range("K22").Value = 638.8148
$ws.Range("L22").Value = 1208.1666
$ws.Range("M22").Value = -288.8148
$ws.Range("N22").Value = -1908.1666
$ws.Range("H31").Value = 54477.184
$ws.Range("I31").Value = 4203.2
$ws.Range("J31").Value = 96372.164
$ws.Range("K31").Value = 4203.2
$ws.Range("L31").Value = 96372.164
$ws.Range("M31").Value = -3908.2
$ws.Range("N31").Value = -96962.164
$ws.Range("H34").Value = 54477.184
$ws.Range("I34").Value = 4203.2
$ws.Range("J34").Value = 96372.164
$ws.Range("K34").Value = 4203.2
$ws.Range("L34").Value = 96372.164
$ws.Range("M34").Value = -4001.2
$ws.Range("N34").Value = -96776.164
$ws.Range("H94").Value = 2467.25
$ws.Range("I94").Value = 1518
$ws.Range("J94").Value = 2941.875
$ws.Range("K94").Value = 1518
$ws.Range("L94").Value = 2941.875
$ws.Range("M94").Value = -1067
$ws.Range("N94").Value = -3843.875
$ws.Range("H129").Value = 51000
$ws.Range("J129").Value = 51000
$ws.Range("L129").Value = 51000
$ws.Range("N129").Value = -61000
$ws.Range("H135").Value = 62499.25
$ws.Range("J135").Value = 69999
$ws.Range("L135").Value = 69999
$ws.Range("N135").Value = -80139

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 19014
$ws.Range("J62").Value = 19014
$ws.Range("L62").Value = 57042
$ws.Range("N62").Value = -58414
$ws.Range("H65").Value = 19014
$ws.Range("J65").Value = 19014
$ws.Range("L65").Value = 171126
$ws.Range("N65").Value = -177990
$ws.Range("H133").Value = 5831.143
$ws.Range("I133").Value = 3631.3333
$ws.Range("K133").Value = 10893.9999
$ws.Range("M133").Value = -5833.999899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1197.3
$ws.Range("I97").Value = 1011
$ws.Range("J97").Value = 1277.1428
$ws.Range("K97").Value = 1011
$ws.Range("L97").Value = 1277.1428
$ws.Range("M97").Value = -515
$ws.Range("N97").Value = -2269.1428
$ws.Range("H113").Value = 6476.4614
$ws.Range("I113").Value = 1839.4
$ws.Range("K113").Value = 1839.4
$ws.Range("M113").Value = 330.5999999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4369.2
$ws.Range("I132").Value = 2532.25
$ws.Range("J132").Value = 7124.625
$ws.Range("K132").Value = 7596.75
$ws.Range("L132").Value = 21373.875
$ws.Range("M132").Value = -5066.75
$ws.Range("N132").Value = -26433.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1555.9524
$ws.Range("I107").Value = 1801.0714
$ws.Range("J107").Value = 1065.7142
$ws.Range("K107").Value = 5403.2142
$ws.Range("L107").Value = 3197.1426
$ws.Range("M107").Value = -3483.2142
$ws.Range("N107").Value = -7037.142599999999
$ws.Range("H132").Value = 4611.1665
$ws.Range("J132").Value = 6663.25
$ws.Range("L132").Value = 19989.75
$ws.Range("N132").Value = -25049.75
$ws.Range("H140").Value = 61099.75
$ws.Range("J140").Value = 61099.75
$ws.Range("L140").Value = 61099.75
$ws.Range("N140").Value = -71459.75
